# Applies the "Diverse small changes as requested" commit to ART_dict.xlsx
# - Renames the German "Literatenquiz" test to "Autor:innenquiz" (various strings)
# - Tweaks a couple of German prompt strings
# - Adds "no/don't know" nuance to the NO answer strings (DE + EN)
# - Adds a new FINISHED row (row 16) with DE/EN strings

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TESTNAME): German test title
$ws.Cells.Item(2, 2).Value = "Autor:innenquiz"

# Row 4 (INSTRUCTIONS_SINGLE): German instructions text rewritten
$ws.Cells.Item(4, 2).Value = "Im Folgenden zeigen wir Ihnen eine Reihe von Personennamen und fragen Sie jeweils, ob es sich bei diesen Personen um Autor:innen von literarischer Prosa, Gedichten oder Dramen handelt. Dies trifft nur für einige der gelisteten Namen zu.<br/>`nBitte klicken Sie nur dann „ja“ an, wenn Sie <strong>sich sicher sind</strong>, dass es sich um die Namen literarischer Autor*innen handelt. Wenn Sie **nicht wissen**, ob es sich um eine/n Autor*in handelt, oder wissen, dass es sich nicht um eine/n Autor*in handelt, klicken Sie „nein/weiß nicht“. Bitte raten Sie nicht. <br/>Sie haben für jede Antwort <strong>maximal 10 Sekunden Zeit</strong>. Wenn Sie sich innerhalb dieser Zeit nicht entschieden haben, wird automatisch der nächste Name angezeigt."

# Row 7 (PROMPT_SINGLE): German prompt wording
$ws.Cells.Item(7, 2).Value = "Ist <b>{{name}}</b> ein:e Autor:in?<br/> Klicken Sie Ja oder Nein, sie haben {{time_out}} Sekunden Zeit zu antworten."

# Row 13 (WELCOME): German welcome title
$ws.Cells.Item(13, 2).Value = "Test: Autor:innenquiz"

# Row 15 (NO): DE + EN "no / don't know" wording
$ws.Cells.Item(15, 2).Value = "Nein/Weiß nicht"
$ws.Cells.Item(15, 3).Value = "No/Don't know"

# New row 16 (FINISHED)
$ws.Cells.Item(16, 1).Value = "FINISHED"
$ws.Cells.Item(16, 2).Value = "Das Autor:innenquize ist nun beendet."
$ws.Cells.Item(16, 3).Value = "You finished the Artist Recognition Test"

# Update the selection to match the final saved state (A16)
$ws.Range("A16").Select()
